# FinalProject/최종 시험-문제참고용.xlsx update
#
# Matches the commit: scoring sheet touch-up after adding Sound, Sword,
# Pause, Record, GunEnemy, PlayerDie and the weapon-cooldown display.
#  - Question #13 (row 28) now has its "student" score filled in (3,
#    matching the base score in column C).
#  - The four "예정" (pending) placeholders in column E for rows 46-49
#    are cleared out now that those items are done.
#  - The running total in D51 recalculates automatically.
#  - Leave the cursor/selection where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Question 13 ("List class 사용"): student score now recorded as 3.
$ws.Range("D28").Value = 3

# Those four rows are no longer "예정" (pending) - clear the markers.
$ws.Range("E46").ClearContents()
$ws.Range("E47").ClearContents()
$ws.Range("E48").ClearContents()
$ws.Range("E49").ClearContents()

# Park the selection/cursor where the author left off.
$ws.Range("N39").Select()
